# Updated cryptos list on Fri Sep  8 05:01:11 UTC 2023 with GitHub Actions
#
# Refreshes the crypto tracker sheet (prices in column D, 1h volume % in
# column E) with the latest scraped figures, and corrects the row-12/13
# ordering for Polkadot vs. Wrapped liquid staked Ether 2.0 (name + link
# swapped back to match the refreshed ranking).
#
# Price-column (D) values are numeric-looking text (e.g. "26.314.05",
# "1.00") that must stay literal strings -- force the Text number format
# before writing so Excel's COM layer doesn't auto-coerce them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.314.05'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.649.01'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.49'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.507'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.258'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.04'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.877.23'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.31'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.664.16'
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.551'
$ws.Range("E15").Value = '  -1.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.286.64'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.999'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '196.84'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.45'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.09'
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.35'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -2.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.11'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.66'
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.26'
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  +1.74%  '
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.916'
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.141.13'
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.556'
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.50'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.67'
$ws.Range("E42").Value = '  +1.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.48'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.785.30'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.53'
$ws.Range("E47").Value = '  +1.89%  '
$ws.Range("E48").Value = '  +2.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0517'
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.75'
$ws.Range("E50").Value = '  +3.05%  '
$ws.Range("E51").Value = '  -0.43%  '
